# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) for the worker's account-statement rows
# (16-20) was re-sorted from descending (2009..2005) to ascending
# (2005..2009).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2005"
$ws.Range("E17").Value = "2006"
$ws.Range("E18").Value = "2007"
$ws.Range("E19").Value = "2008"
$ws.Range("E20").Value = "2009"
